$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 6814327
$ws.Range("E2").Value = "NS Mura"
$ws.Range("F2").Value = "NK Domzale"
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = "A"
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = 3.3
$ws.Range("N2").Value = 3.4
$ws.Range("O2").Value = 1.909
$ws.Range("P2").Value = 3.4
$ws.Range("R2").Value = -0.5
$ws.Range("S2").Value = 1.95
$ws.Range("T2").Value = 1.85
$ws.Range("U2").Value = 2.5
$ws.Range("V2").Value = 1.9
$ws.Range("W2").Value = 1.9
$ws.Range("X2").Value = -1
$ws.Range("Z2").Value = 2.75
$ws.Range("AA2").Value = -1
$ws.Range("AB2").Value = 0.8500000000000001
$ws.Range("AC2").Value = 0.8999999999999999
$ws.Range("AD2").Value = -1

# Row 3
$ws.Range("B3").Value = 6816473
$ws.Range("E3").Value = "NK Bravo"
$ws.Range("F3").Value = "NK Rogaska"
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "H"
$ws.Range("L3").Value = 1.8
$ws.Range("M3").Value = 3.5
$ws.Range("N3").Value = 4
$ws.Range("O3").Value = 2.05
$ws.Range("P3").Value = 3
$ws.Range("R3").Value = -0.25
$ws.Range("S3").Value = 1.75
$ws.Range("T3").Value = 2.05
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.95
$ws.Range("W3").Value = 1.85
$ws.Range("X3").Value = 1.05
$ws.Range("Z3").Value = -1
$ws.Range("AA3").Value = 0.75
$ws.Range("AB3").Value = -1
$ws.Range("AC3").Value = -0.5
$ws.Range("AD3").Value = 0.425

# Row 9
$ws.Range("B9").Value = 6814328
$ws.Range("E9").Value = "NK Domzale"
$ws.Range("F9").Value = "NK Bravo"
$ws.Range("H9").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = "D"
$ws.Range("L9").Value = 2.35
$ws.Range("M9").Value = 3.1
$ws.Range("N9").Value = 2.9
$ws.Range("O9").Value = 2.15
$ws.Range("P9").Value = 3.1
$ws.Range("Q9").Value = 3.3
$ws.Range("R9").Value = -0.25
$ws.Range("S9").Value = 1.925
$ws.Range("T9").Value = 1.875
$ws.Range("U9").Value = 2.25
$ws.Range("V9").Value = 1.95
$ws.Range("W9").Value = 1.85
$ws.Range("X9").Value = -1
$ws.Range("Y9").Value = 2.1
$ws.Range("AB9").Value = 0.4375
$ws.Range("AC9").Value = -0.5
$ws.Range("AD9").Value = 0.425

# Row 10
$ws.Range("B10").Value = 6814330
$ws.Range("E10").Value = "NK Maribor"
$ws.Range("F10").Value = "NK Aluminij"
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = "H"
$ws.Range("L10").Value = 1.363
$ws.Range("M10").Value = 4.5
$ws.Range("N10").Value = 7
$ws.Range("O10").Value = 1.4
$ws.Range("P10").Value = 4.5
$ws.Range("Q10").Value = 7
$ws.Range("R10").Value = -1.25
$ws.Range("S10").Value = 1.85
$ws.Range("T10").Value = 1.95
$ws.Range("U10").Value = 2.75
$ws.Range("V10").Value = 1.8
$ws.Range("W10").Value = 2
$ws.Range("X10").Value = 0.3999999999999999
$ws.Range("Y10").Value = -1
$ws.Range("AB10").Value = 0.475
$ws.Range("AC10").Value = -1
$ws.Range("AD10").Value = 1

# Row 174
$ws.Range("B174").Value = 7124153
$ws.Range("E174").Value = "NK Aluminij"
$ws.Range("F174").Value = "NK Domzale"
$ws.Range("H174").Value = 3
$ws.Range("J174").Value = 3
$ws.Range("L174").Value = 2
$ws.Range("M174").Value = 3.6
$ws.Range("N174").Value = 3
$ws.Range("O174").Value = 1.333
$ws.Range("P174").Value = 4.75
$ws.Range("Q174").Value = 7
$ws.Range("R174").Value = -1.5
$ws.Range("S174").Value = 1.95
$ws.Range("T174").Value = 1.85
$ws.Range("U174").Value = 3.25
$ws.Range("V174").Value = 1.95
$ws.Range("W174").Value = 1.85
$ws.Range("Z174").Value = 6
$ws.Range("AB174").Value = 0.8500000000000001
$ws.Range("AC174").Value = 0.95

# Row 175
$ws.Range("B175").Value = 7133777
$ws.Range("E175").Value = "NK Radomlje"
$ws.Range("F175").Value = "NK Celje"
$ws.Range("H175").Value = 1
$ws.Range("I175").Value = 1
$ws.Range("J175").Value = 0
$ws.Range("K175").Value = "D"
$ws.Range("L175").Value = 3.05
$ws.Range("M175").Value = 3.5
$ws.Range("N175").Value = 2
$ws.Range("O175").Value = 2.9
$ws.Range("P175").Value = 3.6
$ws.Range("Q175").Value = 2.1
$ws.Range("R175").Value = 0.25
$ws.Range("S175").Value = 1.9
$ws.Range("T175").Value = 1.9
$ws.Range("U175").Value = 2.75
$ws.Range("V175").Value = 1.8
$ws.Range("W175").Value = 2
$ws.Range("Y175").Value = 2.6
$ws.Range("Z175").Value = -1
$ws.Range("AA175").Value = 0.45
$ws.Range("AB175").Value = -0.5
$ws.Range("AC175").Value = -1
$ws.Range("AD175").Value = 1

# Row 176
$ws.Range("B176").Value = 7124152
$ws.Range("E176").Value = "NS Mura"
$ws.Range("F176").Value = "NK Rogaska"
$ws.Range("H176").Value = 2
$ws.Range("I176").Value = 0
$ws.Range("J176").Value = 2
$ws.Range("K176").Value = "A"
$ws.Range("L176").Value = 2.45
$ws.Range("M176").Value = 3.4
$ws.Range("N176").Value = 2.45
$ws.Range("O176").Value = 3.8
$ws.Range("Q176").Value = 1.8
$ws.Range("R176").Value = 0.5
$ws.Range("S176").Value = 1.975
$ws.Range("T176").Value = 1.825
$ws.Range("U176").Value = 2.5
$ws.Range("Y176").Value = -1
$ws.Range("Z176").Value = 0.8
$ws.Range("AA176").Value = -1
$ws.Range("AB176").Value = 0.825
$ws.Range("AC176").Value = 0.8
$ws.Range("AD176").Value = -1

# Row 177
$ws.Range("B177").Value = 7128629
$ws.Range("E177").Value = "FC Koper"
$ws.Range("F177").Value = "NK Maribor"
$ws.Range("G177").Value = 1
$ws.Range("I177").Value = 1
$ws.Range("J177").Value = 1
$ws.Range("K177").Value = "D"
$ws.Range("L177").Value = 3.9
$ws.Range("M177").Value = 3.7
$ws.Range("N177").Value = 1.7
$ws.Range("O177").Value = 4.333
$ws.Range("Q177").Value = 1.6
$ws.Range("R177").Value = 1
$ws.Range("S177").Value = 1.775
$ws.Range("T177").Value = 2.025
$ws.Range("V177").Value = 1.925
$ws.Range("W177").Value = 1.875
$ws.Range("Y177").Value = 2.9
$ws.Range("Z177").Value = -1
$ws.Range("AA177").Value = 0.7749999999999999
$ws.Range("AB177").Value = -1
$ws.Range("AD177").Value = 0.875

# Row 178
$ws.Range("B178").Value = 7133776
$ws.Range("E178").Value = "Olimpija Ljubljana"
$ws.Range("F178").Value = "NK Bravo"
$ws.Range("G178").Value = 0
$ws.Range("I178").Value = 0
$ws.Range("J178").Value = 0
$ws.Range("K178").Value = "A"
$ws.Range("L178").Value = 1.533
$ws.Range("M178").Value = 4
$ws.Range("N178").Value = 4.75
$ws.Range("O178").Value = 1.615
$ws.Range("Q178").Value = 4.333
$ws.Range("R178").Value = -0.75
$ws.Range("S178").Value = 1.8
$ws.Range("T178").Value = 2
$ws.Range("V178").Value = 1.975
$ws.Range("W178").Value = 1.825
$ws.Range("Y178").Value = -1
$ws.Range("Z178").Value = 3.333
$ws.Range("AA178").Value = -1
$ws.Range("AB178").Value = 1
$ws.Range("AD178").Value = 0.825
